$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 205 (pushing the existing rows 205-239 down to 207-241).
$ws.Rows.Item(205).Insert()
$ws.Rows.Item(205).Insert()

# New row 205: weekly update for "Primera" quality, dated 2021-11-05 (serial 44505).
$ws.Cells.Item(205,1).Value  = 3
$ws.Cells.Item(205,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(205,3).Value  = "Coquimbo"
$ws.Cells.Item(205,4).Value  = 44505
$ws.Cells.Item(205,5).Value  = 5
$ws.Cells.Item(205,6).Value  = 100112013
$ws.Cells.Item(205,7).Value  = "Alcachofa"
$ws.Cells.Item(205,8).Value  = "Española"
$ws.Cells.Item(205,9).Value  = "Primera"
$ws.Cells.Item(205,10).Value = 10300
$ws.Cells.Item(205,11).Value = 300
$ws.Cells.Item(205,12).Value = 330
$ws.Cells.Item(205,13).Value = 317
$ws.Cells.Item(205,14).Value = "$/unidad"
$ws.Cells.Item(205,15).Value = "Llay Llay"
$ws.Cells.Item(205,16).Value = 317
$ws.Cells.Item(205,17).Value = 1
$ws.Cells.Item(205,18).Value = "Hortaliza"

# New row 206: weekly update for "Segunda" quality, dated 2021-11-05 (serial 44505).
$ws.Cells.Item(206,1).Value  = 3
$ws.Cells.Item(206,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(206,3).Value  = "Coquimbo"
$ws.Cells.Item(206,4).Value  = 44505
$ws.Cells.Item(206,5).Value  = 5
$ws.Cells.Item(206,6).Value  = 100112013
$ws.Cells.Item(206,7).Value  = "Alcachofa"
$ws.Cells.Item(206,8).Value  = "Española"
$ws.Cells.Item(206,9).Value  = "Segunda"
$ws.Cells.Item(206,10).Value = 5900
$ws.Cells.Item(206,11).Value = 200
$ws.Cells.Item(206,12).Value = 200
$ws.Cells.Item(206,13).Value = 200
$ws.Cells.Item(206,14).Value = "$/unidad"
$ws.Cells.Item(206,15).Value = "Llay Llay"
$ws.Cells.Item(206,16).Value = 200
$ws.Cells.Item(206,17).Value = 1
$ws.Cells.Item(206,18).Value = "Hortaliza"
